$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the former A2 cell content (row 2 will be removed since it becomes empty)
$ws.Range("A2").Value = $null

# Update A3 text and add B3 value
$ws.Range("A3").Value = "Oxea"
$ws.Range("B3").Value = 500

$ws.Range("B4").Select() | Out-Null
